$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 328
$ws.Range("D2").Value = 310.5

# Row 3
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = 43.5

# Row 4
$ws.Range("B4").Value = 1673
$ws.Range("C4").Value = 1719
$ws.Range("D4").Value = 1696

# Row 6
$ws.Range("C6").Value = 253
$ws.Range("D6").Value = 238.5

# Row 8
$ws.Range("C8").Value = 110
$ws.Range("D8").Value = 61.5

# Row 9
$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 46.5

# Row 10
$ws.Range("C10").Value = 257
$ws.Range("D10").Value = 238.5

# Row 11
$ws.Range("C11").Value = 111
$ws.Range("D11").Value = 62.5

# Row 12
$ws.Range("C12").Value = 37
$ws.Range("D12").Value = 35.5

# Row 13
$ws.Range("C13").Value = 161
$ws.Range("D13").Value = 147

# Row 14
$ws.Range("C14").Value = 106
$ws.Range("D14").Value = 103.5

# Row 17
$ws.Range("C17").Value = 34
$ws.Range("D17").Value = 34

# Row 18
$ws.Range("C18").Value = 147
$ws.Range("D18").Value = 119.5

# Row 19
$ws.Range("C19").Value = 121
$ws.Range("D19").Value = 118

# Row 20
$ws.Range("C20").Value = 151
$ws.Range("D20").Value = 138.5

# Row 22
$ws.Range("C22").Value = 36
$ws.Range("D22").Value = 23

# Row 23
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 32.5

# Row 24
$ws.Range("C24").Value = 234
$ws.Range("D24").Value = 239.5

# Row 25
$ws.Range("C25").Value = 54
$ws.Range("D25").Value = 50

# Row 27
$ws.Range("C27").Value = 528
$ws.Range("D27").Value = 503

# Row 28
$ws.Range("B28").Value = 169.56
$ws.Range("C28").Value = 184
